# Simplify project to use custom fields:
# Add a new "filter_by" lookup sheet (index/value/column) ahead of the
# existing "time_period" sheet, so "filter_by" becomes the first/active tab.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts a new sheet immediately before the active sheet,
# which is exactly the "filter_by" ahead of "time_period" ordering we need.
$filterBy = $wb.Worksheets.Add()
$filterBy.Name = "filter_by"

# Header row.
$filterBy.Range("A1").Value = "index"
$filterBy.Range("B1").Value = "value"
$filterBy.Range("C1").Value = "column"

# Column A (index) first, then column B (value), then column C (column) —
# matches the shared-string interning order of the authored workbook.
$filterBy.Range("A2").Value = 0
$filterBy.Range("A3").Value = 1
$filterBy.Range("A4").Value = 2

$filterBy.Range("B2").Value = "year"
$filterBy.Range("B3").Value = "quarter"
$filterBy.Range("B4").Value = "month"

$filterBy.Range("C2").Value = "start_of_year"
$filterBy.Range("C3").Value = "start_of_quarter"
$filterBy.Range("C4").Value = "start_of_month"
